# Update "想去人数" (interested-count) figures across the workbook.
# Sheet "全部类型" (All Types) mirrors rows from the other three sheets, so
# the same events get updated there too (with their own independent counts).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1919
$ws1.Range("F12").Value = 1216
$ws1.Range("F15").Value = 2292
$ws1.Range("F17").Value = 847
$ws1.Range("F18").Value = 1049
$ws1.Range("F23").Value = 127
$ws1.Range("F25").Value = 1289
$ws1.Range("F29").Value = 126
$ws1.Range("F40").Value = 2162
$ws1.Range("F43").Value = 1941

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 0
$ws2.Range("F31").Value = 31

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value  = 4788
$ws3.Range("F9").Value  = 634
$ws3.Range("F12").Value = 530
$ws3.Range("F13").Value = 1246
$ws3.Range("F15").Value = 975

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1919
$ws4.Range("F5").Value  = 4788
$ws4.Range("F6").Value  = 634
$ws4.Range("F10").Value = 530
$ws4.Range("F11").Value = 1246
$ws4.Range("F19").Value = 1216
$ws4.Range("F21").Value = 975
$ws4.Range("F22").Value = 975
$ws4.Range("F23").Value = 2292
$ws4.Range("F26").Value = 847
$ws4.Range("F27").Value = 1049
$ws4.Range("F32").Value = 127
$ws4.Range("F33").Value = 1289
$ws4.Range("F36").Value = 126
$ws4.Range("F46").Value = 2162
$ws4.Range("F49").Value = 1942
